$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the student ID ("学号") values in column A (rows 2-9).
# A leading apostrophe keeps the numeric-looking string stored as text
# (preserving the leading zeros), matching the existing text formatting
# of these cells.
$ws.Range("A2").Value = "'00100012"
$ws.Range("A3").Value = "'00100013"
$ws.Range("A4").Value = "'00100014"
$ws.Range("A5").Value = "'00100015"
$ws.Range("A6").Value = "'00100016"
$ws.Range("A7").Value = "'00100017"
$ws.Range("A8").Value = "'00100018"
$ws.Range("A9").Value = "'00100019"

# Update the sheet view's selection to A2:A9 with A2 as the active cell.
$ws.Range("A2:A9").Select()
